$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the paragraph that currently reads "Example Formatted Output ..."
# (it is the one that is being split into three paragraphs by this edit).
# ------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Example Formatted Output*") {
        $target = $p
        break
    }
}

$r = $target.Range

# ------------------------------------------------------------------
# Replace that single paragraph with the three paragraphs described by the
# diff:
#   1) " " (a lone space, carrying the <w:lastRenderedPageBreak/> that used
#      to sit on the "Example Formatted Output" run) followed by the new
#      bold "Compare all Lists" heading and its description run.
#   2) a blank spacer paragraph
#   3) the original "Example Formatted Output ..." paragraph, now starting
#      with a plain (no lastRenderedPageBreak) bold run.
# The bookmarks around "Compare all Lists...multivector." are added
# afterwards via Bookmarks.Add so the engine assigns/renumbers their
# w:id values the same way Word itself would.
# ------------------------------------------------------------------
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p>' +
'<w:pPr><w:ind w:left="360"/></w:pPr>' +
'<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve"> </w:t></w:r>' +
'<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Compare all Lists</w:t></w:r>' +
'<w:r><w:t>' + [char]0x2013 + 'A Notebook to supplement the Working Tutorial that displays, side-by-side, examples of the different lists that this GA package can generate from a given multivector.</w:t></w:r>' +
'</w:p>' +
'<w:p>' +
'<w:pPr><w:ind w:left="360"/></w:pPr>' +
'</w:p>' +
'<w:p>' +
'<w:pPr><w:ind w:left="360"/></w:pPr>' +
'<w:r><w:rPr><w:b/></w:rPr><w:t>Example Formatted Output</w:t></w:r>' +
'<w:r><w:t xml:space="preserve"> - Shows how function </w:t></w:r>' +
'<w:proofErr w:type="spellStart"/>' +
'<w:r><w:t>AtomCoefG</w:t></w:r>' +
'<w:proofErr w:type="spellEnd"/>' +
'<w:r><w:t xml:space="preserve"> can be used to format output </w:t></w:r>' +
'<w:r><w:t>to line up in</w:t></w:r>' +
'<w:r><w:t xml:space="preserve"> a table format</w:t></w:r>' +
'</w:p>' +
'</w:body>' +
'</w:document>' +
'</pkg:xmlData>' +
'</pkg:part>' +
'</pkg:package>'

$r.InsertXML($xml)

# ------------------------------------------------------------------
# Now wrap the freshly-inserted "Compare all Lists ... multivector." text
# with the two new bookmarks (OLE_LINK12 / OLE_LINK13). Adding them through
# Bookmarks.Add lets the document engine pick free w:id values and bump any
# existing bookmarks that collide, exactly like Word's own bookmark-insert
# behaviour.
# ------------------------------------------------------------------
$bmStart = $null
$bmEnd = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Compare all Lists*multivector.*") {
        $t = $p.Range.Text
        $relStart = $t.IndexOf("Compare all Lists")
        $relEnd = $t.IndexOf("multivector.") + "multivector.".Length
        $bmStart = $p.Range.Start + $relStart
        $bmEnd = $p.Range.Start + $relEnd
        break
    }
}

$bmRange = $d.Range($bmStart, $bmEnd)
Write-Output ("Bookmark range text: [" + $bmRange.Text + "]")

$d.Bookmarks.Add("OLE_LINK12", $bmRange)
$d.Bookmarks.Add("OLE_LINK13", $bmRange)
